# Insert a new row above row 1, shifting all existing rows (1-25) down to (2-26).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert()

# Copy the (now shifted-down) former header row's formatting onto the brand-new
# row 1 so it keeps the bold/bordered "header" look (style index 1), then
# overwrite the old header row's style back to the default "Normal" style.
$ws.Range("A2:K2").Copy()
$ws.Range("A1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A2:K2").Style = "Normal"

# Fill the new row 1 with a simple numeric column index (0-based).
for ($c = 1; $c -le 11; $c++) {
    $ws.Cells.Item(1, $c).Value = $c - 1
}

# The old header row (now row 2) no longer carries the "thread_size" /
# "material_surface" labels that used to live in J1/K1.
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
